# Runnning all the test cases
# Update the "Runmode" column (C) on the "Test Suite" sheet from "N" to "Y"
# for all rows except row 3 (already "Y"), and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("C2").Value = "Y"
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "Y"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

$ws.Range("C2:C7").Select()
